# Update FFXIV leve/market profit figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# to reflect the latest scraped Universalis market prices (scheduled runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 140092.64
$ws.Range("I40").Value = 1002429.7
$ws.Range("J40").Value = 3934.158
$ws.Range("K40").Value = 1002429.7
$ws.Range("L40").Value = 3934.158
$ws.Range("M40").Value = -1002254.7
$ws.Range("N40").Value = -4284.157999999999
$ws.Range("H62").Value = 6807.0527
$ws.Range("I62").Value = 5336.3335
$ws.Range("K62").Value = 5336.3335
$ws.Range("M62").Value = -4712.3335
$ws.Range("H65").Value = 6807.0527
$ws.Range("I65").Value = 5336.3335
$ws.Range("K65").Value = 26681.6675
$ws.Range("M65").Value = -23561.6675
$ws.Range("H111").Value = 1913.7142
$ws.Range("I111").Value = 1952.8
$ws.Range("J111").Value = 1816.0
$ws.Range("K111").Value = 5858.4
$ws.Range("L111").Value = 5448.0
$ws.Range("M111").Value = -2791.4
$ws.Range("N111").Value = -11582.0
$ws.Range("H116").Value = 13196.5625
$ws.Range("I116").Value = 15597.917
$ws.Range("K116").Value = 15597.917
$ws.Range("M116").Value = -12155.917
$ws.Range("H118").Value = 1674.4
$ws.Range("I118").Value = 998.0
$ws.Range("J118").Value = 3252.6667
$ws.Range("K118").Value = 2994.0
$ws.Range("L118").Value = 9758.000100000001
$ws.Range("M118").Value = -1337.0
$ws.Range("N118").Value = -13072.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 479.85715
$ws.Range("I5").Value = 479.85715
$ws.Range("J5").Value = 0.0
$ws.Range("K5").Value = 479.85715
$ws.Range("L5").Value = 0.0
$ws.Range("M5").Value = -367.85715
$ws.Range("N5").Value = ""
$ws.Range("H11").Value = 0.0
$ws.Range("I11").Value = 0.0
$ws.Range("K11").Value = 0.0
$ws.Range("M11").Value = ""
$ws.Range("H61").Value = 2390.0
$ws.Range("I61").Value = 2207.2778
$ws.Range("K61").Value = 2207.2778
$ws.Range("M61").Value = -1995.2778
$ws.Range("H74").Value = 3170.4
$ws.Range("I74").Value = 2504.75
$ws.Range("K74").Value = 2504.75
$ws.Range("M74").Value = -1630.75
$ws.Range("H77").Value = 3170.4
$ws.Range("I77").Value = 2504.75
$ws.Range("K77").Value = 12523.75
$ws.Range("M77").Value = -8155.75
$ws.Range("H136").Value = 2390.0
$ws.Range("I136").Value = 2207.2778
$ws.Range("K136").Value = 6621.8334
$ws.Range("M136").Value = -4071.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 479.85715
$ws.Range("I4").Value = 479.85715
$ws.Range("J4").Value = 0.0
$ws.Range("K4").Value = 479.85715
$ws.Range("L4").Value = 0.0
$ws.Range("M4").Value = -364.85715
$ws.Range("N4").Value = ""
$ws.Range("H12").Value = 0.0
$ws.Range("I12").Value = 0.0
$ws.Range("K12").Value = 0.0
$ws.Range("M12").Value = ""
$ws.Range("H122").Value = 53330.0
$ws.Range("J122").Value = 29995.0
$ws.Range("L122").Value = 29995.0
$ws.Range("N122").Value = -39795.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0.0
$ws.Range("I62").Value = 0.0
$ws.Range("K62").Value = 0.0
$ws.Range("M62").Value = ""
$ws.Range("H65").Value = 0.0
$ws.Range("I65").Value = 0.0
$ws.Range("K65").Value = 0.0
$ws.Range("M65").Value = ""
$ws.Range("H86").Value = 5751.75
$ws.Range("I86").Value = 6503.5
$ws.Range("K86").Value = 6503.5
$ws.Range("M86").Value = -5380.5
$ws.Range("H89").Value = 5751.75
$ws.Range("I89").Value = 6503.5
$ws.Range("K89").Value = 32517.5
$ws.Range("M89").Value = -26901.5
$ws.Range("H99").Value = 3000.0
$ws.Range("I99").Value = 0.0
$ws.Range("K99").Value = 0.0
$ws.Range("M99").Value = ""
$ws.Range("H126").Value = 3000.0
$ws.Range("I126").Value = 0.0
$ws.Range("K126").Value = 0.0
$ws.Range("M126").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 161.25
$ws.Range("I7").Value = 0.0
$ws.Range("J7").Value = 161.25
$ws.Range("K7").Value = 0.0
$ws.Range("L7").Value = 483.75
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = -707.75
$ws.Range("H13").Value = 0.0
$ws.Range("I13").Value = 0.0
$ws.Range("J13").Value = 0.0
$ws.Range("K13").Value = 0.0
$ws.Range("L13").Value = 0.0
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("H50").Value = 2929.1428
$ws.Range("I50").Value = 2452.5
$ws.Range("J50").Value = 3119.8
$ws.Range("K50").Value = 7357.5
$ws.Range("L50").Value = 9359.400000000001
$ws.Range("M50").Value = -6876.5
$ws.Range("N50").Value = -10321.4
$ws.Range("H53").Value = 2929.1428
$ws.Range("I53").Value = 2452.5
$ws.Range("J53").Value = 3119.8
$ws.Range("K53").Value = 7357.5
$ws.Range("L53").Value = 9359.400000000001
$ws.Range("M53").Value = -6876.5
$ws.Range("N53").Value = -10321.4
$ws.Range("H109").Value = 0.0
$ws.Range("I109").Value = 0.0
$ws.Range("J109").Value = 0.0
$ws.Range("K109").Value = 0.0
$ws.Range("L109").Value = 0.0
$ws.Range("M109").Value = ""
$ws.Range("N109").Value = ""
$ws.Range("H140").Value = 13335152.0
$ws.Range("I140").Value = 13335152.0
$ws.Range("K140").Value = 40005456.0
$ws.Range("M140").Value = -40000276.0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 520.0
$ws.Range("I107").Value = 617.2857
$ws.Range("J107").Value = 349.75
$ws.Range("K107").Value = 617.2857
$ws.Range("L107").Value = 349.75
$ws.Range("M107").Value = 1302.7143
$ws.Range("N107").Value = -4189.75
$ws.Range("H126").Value = 4632.304
$ws.Range("I126").Value = 4708.4443
$ws.Range("K126").Value = 14125.3329
$ws.Range("M126").Value = -11655.3329
$ws.Range("H132").Value = 2215.72
$ws.Range("I132").Value = 2215.72
$ws.Range("J132").Value = 0.0
$ws.Range("K132").Value = 6647.16
$ws.Range("L132").Value = 0.0
$ws.Range("M132").Value = -4117.16
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2557.5
$ws.Range("I16").Value = 2620.5
$ws.Range("J16").Value = 2400.0
$ws.Range("K16").Value = 2620.5
$ws.Range("L16").Value = 2400.0
$ws.Range("M16").Value = -2450.5
$ws.Range("N16").Value = -2740.0
$ws.Range("H22").Value = 3500.0
$ws.Range("J22").Value = 2000.0
$ws.Range("L22").Value = 2000.0
$ws.Range("N22").Value = -2590.0
$ws.Range("H27").Value = 3500.0
$ws.Range("J27").Value = 2000.0
$ws.Range("L27").Value = 2000.0
$ws.Range("N27").Value = -2214.0
$ws.Range("H40").Value = 6677.4
$ws.Range("I40").Value = 6523.074
$ws.Range("J40").Value = 8066.3335
$ws.Range("K40").Value = 6523.074
$ws.Range("L40").Value = 8066.3335
$ws.Range("M40").Value = -6387.074
$ws.Range("N40").Value = -8338.3335
$ws.Range("H122").Value = 17079.309
$ws.Range("I122").Value = 17079.309
$ws.Range("K122").Value = 51237.927
$ws.Range("M122").Value = -48787.927
$ws.Range("H132").Value = 592052.9
$ws.Range("I132").Value = 717657.06
$ws.Range("K132").Value = 2152971.18
$ws.Range("M132").Value = -2150441.18
$ws.Range("H136").Value = 5101.643
$ws.Range("I136").Value = 3986.25
$ws.Range("K136").Value = 11958.75
$ws.Range("M136").Value = -9408.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 0.0
$ws.Range("J99").Value = 0.0
$ws.Range("L99").Value = 0.0
$ws.Range("N99").Value = ""
$ws.Range("H107").Value = 456.75
$ws.Range("I107").Value = 364.33334
$ws.Range("J107").Value = 575.5714
$ws.Range("K107").Value = 1093.00002
$ws.Range("L107").Value = 1726.7142
$ws.Range("M107").Value = 826.99998
$ws.Range("N107").Value = -5566.7142
$ws.Range("H128").Value = 50078000.0
$ws.Range("J128").Value = 50078000.0
$ws.Range("L128").Value = 50078000.0
$ws.Range("N128").Value = -50087960.0
$ws.Range("H136").Value = 26507.928
$ws.Range("I136").Value = 1942.9048
$ws.Range("J136").Value = 51072.953
$ws.Range("K136").Value = 5828.7144
$ws.Range("L136").Value = 153218.859
$ws.Range("M136").Value = -3278.7144
$ws.Range("N136").Value = -158318.859
